$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Shortlist")

$ws.Range("A2").Value = "Belegschaft Lieferkette"
$ws.Range("B2").Value = "Arbeitsbedingungen"
$ws.Range("C2").Value = "Angemessene Entlohnung"

$ws.Range("A3").Value = "Biodiversität"
$ws.Range("B3").Value = "Direkte Ursachen des Biodiversitätsverlusts"
$ws.Range("C3").Value = "Direkte Ausbeutung"

$ws.Range("A4").Value = "Externes Thema 1"
$ws.Range("B4").Value = "Externes Unterthema 1"
$ws.Range("C4").Value = "Externes Unter-Unterthema 1"
